$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet from "Sheet1" to "Table" (updates the _FilterDatabase
# defined name reference automatically since it tracks the sheet).
$ws.Name = "Table"

# Update the active selection on the sheet from D21 to E12.
$ws.Range("E12").Select()
